$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values that look numeric (e.g. "4.81") stay as text,
# matching the original inline-string cell type used throughout the sheet.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '68.603.37'
$ws.Range("E2").Value = '  -1.55%  '

$ws.Range("D3").Value = '2.456.39'
$ws.Range("E3").Value = '  -2.13%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '564.35'
$ws.Range("E5").Value = '  -1.94%  '

$ws.Range("D6").Value = '162.96'
$ws.Range("E6").Value = '  -2.27%  '

$ws.Range("D8").Value = '0.506'
$ws.Range("E8").Value = '  -1.59%  '

$ws.Range("E9").Value = '  -6.77%  '

$ws.Range("E10").Value = '  -1.90%  '

$ws.Range("E11").Value = '  -4.36%  '

$ws.Range("D12").Value = '4.81'
$ws.Range("E12").Value = '  -2.60%  '

$ws.Range("D13").Value = '2.903.70'
$ws.Range("E13").Value = '  -2.23%  '

$ws.Range("D14").Value = '68.414.50'
$ws.Range("E14").Value = '  -1.65%  '

$ws.Range("E15").Value = '  -3.92%  '

$ws.Range("D16").Value = '23.69'
$ws.Range("E16").Value = '  -4.76%  '

$ws.Range("D17").Value = '2.463.17'
$ws.Range("E17").Value = '  -2.26%  '

$ws.Range("D18").Value = '10.99'
$ws.Range("E18").Value = '  -2.24%  '

$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = '344.21'
$ws.Range("E19").Value = '  -1.69%  '

$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '7.19'
$ws.Range("E20").Value = '  -4.40%  '

$ws.Range("D21").Value = '3.82'
$ws.Range("E21").Value = '  -2.56%  '

$ws.Range("E22").Value = '  -4.30%  '

$ws.Range("E23").Value = '  +0.01%  '

$ws.Range("D24").Value = '68.17'
$ws.Range("E24").Value = '  -2.69%  '

$ws.Range("E25").Value = '  -4.81%  '

$ws.Range("D26").Value = '1.05'
$ws.Range("E26").Value = '  +4.77%  '

$ws.Range("D27").Value = '2.582.87'
$ws.Range("E27").Value = '  -2.64%  '

$ws.Range("D28").Value = '8.23'
$ws.Range("E28").Value = '  -6.97%  '

$ws.Range("D29").Value = '0.0₃0843'
$ws.Range("E29").Value = '  -5.64%  '

$ws.Range("D30").Value = '7.32'
$ws.Range("E30").Value = '  -7.00%  '

$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").Value = '437.47'
$ws.Range("E31").Value = '  -4.79%  '

$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").Value = '1.19'
$ws.Range("E32").Value = '  -2.84%  '

$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.05%  '

$ws.Range("D34").Value = '1.69'
$ws.Range("E34").Value = '  -2.85%  '

$ws.Range("E35").Value = '  +97.86%  '

$ws.Range("D36").Value = '156.79'
$ws.Range("E36").Value = '  -2.00%  '

$ws.Range("E37").Value = '  -0.37%  '

$ws.Range("E38").Value = '  +0.03%  '

$ws.Range("D39").Value = '0.110'
$ws.Range("E39").Value = '  -5.69%  '

$ws.Range("E40").Value = '  -3.25%  '

$ws.Range("E41").Value = '  -3.69%  '

$ws.Range("D42").Value = '4.51'

$ws.Range("D43").Value = '1.54'
$ws.Range("E43").Value = '  -3.85%  '

$ws.Range("D44").Value = '1.11'
$ws.Range("E44").Value = '  +1.78%  '

$ws.Range("D45").Value = '2.10'
$ws.Range("E45").Value = '  -4.96%  '

$ws.Range("D46").Value = '135.31'
$ws.Range("E46").Value = '  -4.85%  '

$ws.Range("D47").Value = '3.38'
$ws.Range("E47").Value = '  -2.72%  '

$ws.Range("E48").Value = '  -5.70%  '

$ws.Range("E49").Value = '  -2.09%  '

$ws.Range("D50").Value = '0.563'
$ws.Range("E50").Value = '  -2.82%  '

$ws.Range("D51").Value = '0.0915'
$ws.Range("E51").Value = '  -1.52%  '
